$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Slit3"
$ws.Cells.Item(2, 3).Value = "Robo2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 3.590118666666667
$ws.Cells.Item(2, 8).Value = 10.770356
$ws.Cells.Item(2, 9).Value = 0.03935865391742773
$ws.Cells.Item(2, 10).Value = 0.03935865391742773
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.655628666666667
$ws.Cells.Item(2, 14).Value = 4.966886000000001
$ws.Cells.Item(2, 15).Value = 0.9788840386998182
$ws.Cells.Item(2, 16).Value = 0.9788840386998183
$ws.Cells.Item(2, 17).Value = 5.943903381268445
$ws.Cells.Item(2, 18).Value = 53.495130431416
$ws.Cells.Item(2, 19).Value = 0.03852755810448008
$ws.Cells.Item(2, 20).Value = 0.03852755810448008

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Slit3"
$ws.Cells.Item(3, 3).Value = "Robo2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 3.590118666666667
$ws.Cells.Item(3, 8).Value = 10.770356
$ws.Cells.Item(3, 9).Value = 0.03935865391742773
$ws.Cells.Item(3, 10).Value = 0.03935865391742773
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.01121666666666667
$ws.Cells.Item(3, 14).Value = 0.03365
$ws.Cells.Item(3, 15).Value = 0.006631810736596105
$ws.Cells.Item(3, 16).Value = 0.006631810736596105
$ws.Cells.Item(3, 17).Value = 0.04026916437777778
$ws.Cells.Item(3, 18).Value = 0.3624224794
$ws.Cells.Item(3, 19).Value = 0.0002610191436275676
$ws.Cells.Item(3, 20).Value = 0.0002610191436275676

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Slit3"
$ws.Cells.Item(4, 3).Value = "Robo2"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 3.590118666666667
$ws.Cells.Item(4, 8).Value = 10.770356
$ws.Cells.Item(4, 9).Value = 0.03935865391742773
$ws.Cells.Item(4, 10).Value = 0.03935865391742773
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.02449766666666667
$ws.Cells.Item(4, 14).Value = 0.073493
$ws.Cells.Item(4, 15).Value = 0.01448415056358566
$ws.Cells.Item(4, 16).Value = 0.01448415056358566
$ws.Cells.Item(4, 17).Value = 0.08794953038977778
$ws.Cells.Item(4, 18).Value = 0.791545773508
$ws.Cells.Item(4, 19).Value = 0.0005700766693200839
$ws.Cells.Item(4, 20).Value = 0.0005700766693200839

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Slit3"
$ws.Cells.Item(5, 3).Value = "Robo2"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 76.92488366666667
$ws.Cells.Item(5, 8).Value = 230.774651
$ws.Cells.Item(5, 9).Value = 0.8433314202078528
$ws.Cells.Item(5, 10).Value = 0.8433314202078527
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 1.655628666666667
$ws.Cells.Item(5, 14).Value = 4.966886000000001
$ws.Cells.Item(5, 15).Value = 0.9788840386998182
$ws.Cells.Item(5, 16).Value = 0.9788840386998183
$ws.Cells.Item(5, 17).Value = 127.3590425785318
$ws.Cells.Item(5, 18).Value = 1146.231383206786
$ws.Cells.Item(5, 19).Value = 0.8255236665755165
$ws.Cells.Item(5, 20).Value = 0.8255236665755165

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Slit3"
$ws.Cells.Item(6, 3).Value = "Robo2"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 76.92488366666667
$ws.Cells.Item(6, 8).Value = 230.774651
$ws.Cells.Item(6, 9).Value = 0.8433314202078528
$ws.Cells.Item(6, 10).Value = 0.8433314202078527
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.01121666666666667
$ws.Cells.Item(6, 14).Value = 0.03365
$ws.Cells.Item(6, 15).Value = 0.006631810736596105
$ws.Cells.Item(6, 16).Value = 0.006631810736596105
$ws.Cells.Item(6, 17).Value = 0.8628407784611112
$ws.Cells.Item(6, 18).Value = 7.76556700615
$ws.Cells.Item(6, 19).Value = 0.005592814367043279
$ws.Cells.Item(6, 20).Value = 0.005592814367043278

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Slit3"
$ws.Cells.Item(7, 3).Value = "Robo2"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 76.92488366666667
$ws.Cells.Item(7, 8).Value = 230.774651
$ws.Cells.Item(7, 9).Value = 0.8433314202078528
$ws.Cells.Item(7, 10).Value = 0.8433314202078527
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.02449766666666667
$ws.Cells.Item(7, 14).Value = 0.073493
$ws.Cells.Item(7, 15).Value = 0.01448415056358566
$ws.Cells.Item(7, 16).Value = 0.01448415056358566
$ws.Cells.Item(7, 17).Value = 1.884480158438111
$ws.Cells.Item(7, 18).Value = 16.960321425943
$ws.Cells.Item(7, 19).Value = 0.01221493926529307
$ws.Cells.Item(7, 20).Value = 0.01221493926529307

# Row 8
$ws.Cells.Item(8, 1).Value = "M2"
$ws.Cells.Item(8, 2).Value = "Slit3"
$ws.Cells.Item(8, 3).Value = "Robo2"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.041643
$ws.Cells.Item(8, 8).Value = 0.124929
$ws.Cells.Item(8, 9).Value = 0.0004565343313861054
$ws.Cells.Item(8, 10).Value = 0.0004565343313861054
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.655628666666667
$ws.Cells.Item(8, 14).Value = 4.966886000000001
$ws.Cells.Item(8, 15).Value = 0.9788840386998182
$ws.Cells.Item(8, 16).Value = 0.9788840386998183
$ws.Cells.Item(8, 17).Value = 0.068945344566
$ws.Cells.Item(8, 18).Value = 0.6205081010940001
$ws.Cells.Item(8, 19).Value = 0.0004468941701123521
$ws.Cells.Item(8, 20).Value = 0.0004468941701123521

# Row 9
$ws.Cells.Item(9, 1).Value = "M2"
$ws.Cells.Item(9, 2).Value = "Slit3"
$ws.Cells.Item(9, 3).Value = "Robo2"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.041643
$ws.Cells.Item(9, 8).Value = 0.124929
$ws.Cells.Item(9, 9).Value = 0.0004565343313861054
$ws.Cells.Item(9, 10).Value = 0.0004565343313861054
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.01121666666666667
$ws.Cells.Item(9, 14).Value = 0.03365
$ws.Cells.Item(9, 15).Value = 0.006631810736596105
$ws.Cells.Item(9, 16).Value = 0.006631810736596105
$ws.Cells.Item(9, 17).Value = 0.00046709565
$ws.Cells.Item(9, 18).Value = 0.004203860849999999
$ws.Cells.Item(9, 19).Value = 0.000003027649280511098
$ws.Cells.Item(9, 20).Value = 0.000003027649280511098

# Row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Slit3"
$ws.Cells.Item(10, 3).Value = "Robo2"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.041643
$ws.Cells.Item(10, 8).Value = 0.124929
$ws.Cells.Item(10, 9).Value = 0.0004565343313861054
$ws.Cells.Item(10, 10).Value = 0.0004565343313861054
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.02449766666666667
$ws.Cells.Item(10, 14).Value = 0.073493
$ws.Cells.Item(10, 15).Value = 0.01448415056358566
$ws.Cells.Item(10, 16).Value = 0.01448415056358566
$ws.Cells.Item(10, 17).Value = 0.001020156333
$ws.Cells.Item(10, 18).Value = 0.009181406997
$ws.Cells.Item(10, 19).Value = 0.000006612511993242262
$ws.Cells.Item(10, 20).Value = 0.000006612511993242263

# Row 11
$ws.Cells.Item(11, 1).Value = "sCs"
$ws.Cells.Item(11, 2).Value = "Slit3"
$ws.Cells.Item(11, 3).Value = "Robo2"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 10.65883866666667
$ws.Cells.Item(11, 8).Value = 31.976516
$ws.Cells.Item(11, 9).Value = 0.1168533915433334
$ws.Cells.Item(11, 10).Value = 0.1168533915433334
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 1.655628666666667
$ws.Cells.Item(11, 14).Value = 4.966886000000001
$ws.Cells.Item(11, 15).Value = 0.9788840386998182
$ws.Cells.Item(11, 16).Value = 0.9788840386998183
$ws.Cells.Item(11, 17).Value = 17.64707884990845
$ws.Cells.Item(11, 18).Value = 158.823709649176
$ws.Cells.Item(11, 19).Value = 0.1143859198497094
$ws.Cells.Item(11, 20).Value = 0.1143859198497094

# Row 12
$ws.Cells.Item(12, 1).Value = "sCs"
$ws.Cells.Item(12, 2).Value = "Slit3"
$ws.Cells.Item(12, 3).Value = "Robo2"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 10.65883866666667
$ws.Cells.Item(12, 8).Value = 31.976516
$ws.Cells.Item(12, 9).Value = 0.1168533915433334
$ws.Cells.Item(12, 10).Value = 0.1168533915433334
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.01121666666666667
$ws.Cells.Item(12, 14).Value = 0.03365
$ws.Cells.Item(12, 15).Value = 0.006631810736596105
$ws.Cells.Item(12, 16).Value = 0.006631810736596105
$ws.Cells.Item(12, 17).Value = 0.1195566403777778
$ws.Cells.Item(12, 18).Value = 1.0760097634
$ws.Cells.Item(12, 19).Value = 0.0007749495766447471
$ws.Cells.Item(12, 20).Value = 0.0007749495766447471

# Row 13
$ws.Cells.Item(13, 1).Value = "sCs"
$ws.Cells.Item(13, 2).Value = "Slit3"
$ws.Cells.Item(13, 3).Value = "Robo2"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 10.65883866666667
$ws.Cells.Item(13, 8).Value = 31.976516
$ws.Cells.Item(13, 9).Value = 0.1168533915433334
$ws.Cells.Item(13, 10).Value = 0.1168533915433334
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.02449766666666667
$ws.Cells.Item(13, 14).Value = 0.073493
$ws.Cells.Item(13, 15).Value = 0.01448415056358566
$ws.Cells.Item(13, 16).Value = 0.01448415056358566
$ws.Cells.Item(13, 17).Value = 0.2611166767097778
$ws.Cells.Item(13, 18).Value = 2.350050090388
$ws.Cells.Item(13, 19).Value = 0.001692522116979269
$ws.Cells.Item(13, 20).Value = 0.001692522116979269
